# Cambio de orden diapositivas
# 1) Move the slide currently at position 11 (the "phone" screenshot slide,
#    sldId 323) down to position 14, so it now follows the three slides that
#    used to sit right after it (sldId 324, 325, 326). This is a pure
#    reorder - no slide content is added, removed or otherwise modified.
# 2) Refresh the cached "datetimeFigureOut" footer field text (the date
#    auto-field shown on the Slide Master, every slide Layout, and the Notes
#    Master) from 28/7/2025 to 4/8/2025, which is what PowerPoint does to
#    these cached field values whenever the deck is re-saved on a later day.

$p = $ppt.ActivePresentation

# --- 1) Reorder: move slide 11 to slide position 14 ---------------------
$movedSlide = $p.Slides.Item(11)
$movedSlide.MoveTo(14)

# --- 2) Refresh cached date field text -----------------------------------
$dateText = "4/8/2025"
$ppPlaceholderDate = 16

# Slide Master
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
        $sh.TextFrame.TextRange.Text = $dateText
    }
}

# Every slide Layout that belongs to the master
for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
            $sh.TextFrame.TextRange.Text = $dateText
        }
    }
}

# Notes Master
$notesMaster = $p.NotesMaster
for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
    $sh = $notesMaster.Shapes.Item($i)
    if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
        $sh.TextFrame.TextRange.Text = $dateText
    }
}
